$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Mallard (청둥오리) row - it was row 3
$ws.Rows.Item(3).Delete()

# Add new columns F (점수 / Score) and G (누적 확률 / Cumulative probability)
$ws.Range("F1").Value = "점수"
$ws.Range("G1").Value = "누적 확률"
$ws.Range("F1:G1").Style = $ws.Range("E1").Style

$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 0.15

$ws.Range("F3").Value = 30
$ws.Range("G3").Value = 0.03

$ws.Range("F4").Value = 50
$ws.Range("G4").Value = 0.02

$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 0.4

$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 0.1

$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 0.3

$ws.Range("H2").Select()
